$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename headers (nomenclature update: sentence -> logical_block) ---
$ws.Range("E1").Value = "logical_block_uniqueness_ratio"
$ws.Range("J1").Value = "norm(logical_block_uniqueness_ratio)"

# --- Update existing row 2 (goodrelations--v1_type=parsed.ttl) ---
$ws.Range("H2").Value = 0.5705319720580332
$ws.Range("K2").Value = 0.5198498702541686
$ws.Range("L2").Value = 0.8845690472280845
$ws.Range("M2").Value = 0.814017582221236

# --- Update existing row 3 (obo--iao--pno--owl_type=parsed.ttl) ---
$ws.Range("H3").Value = 0.517281553398058
$ws.Range("I3").Value = 0.8076923076923075
$ws.Range("M3").Value = 0.7263652440994914

# --- Insert a new row 4 for resource--leak_type=parsed.ttl, pushing the
#     existing row 4 (rvl_type=parsed.ttl) down to row 5 ---
$ws.Rows(4).Insert()

$ws.Range("A4").Value = "resource--leak_type=parsed.ttl"
$ws.Range("B4").Value = 12072
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3580246913580247
$ws.Range("G4").Value = 34.83363761275684
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0

# --- Update the shifted row 5 (rvl_type=parsed.ttl) values ---
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.8461538461538459
$ws.Range("K5").Value = 0.828153093249641
$ws.Range("L5").Value = 0.9116680167136849
$ws.Range("M5").Value = 0.807827298459074
